{"js": "// Replace the 25 three-digit-by-one-digit multiplication equations in the\n// table with their updated values, preserving all run/paragraph formatting.\nconst replacements = [\n  [\"364\u00d77=2548\", \"334\u00d76=2004\"],\n  [\"171\u00d73=513\", \"291\u00d76=1746\"],\n  [\"748\u00d79=6732\", \"239\u00d73=717\"],\n  [\"716\u00d77=5012\", \"945\u00d77=6615\"],\n  [\"318\u00d75=1590\", \"376\u00d77=2632\"],\n  [\"169\u00d72=338\", \"831\u00d72=1662\"],\n  [\"675\u00d75=3375\", \"713\u00d76=4278\"],\n  [\"201\u00d77=1407\", \"239\u00d73=717\"],\n  [\"568\u00d75=2840\", \"536\u00d76=3216\"],\n  [\"746\u00d72=1492\", \"464\u00d77=3248\"],\n  [\"520\u00d78=4160\", \"502\u00d77=3514\"],\n  [\"588\u00d78=4704\", \"401\u00d78=3208\"],\n  [\"699\u00d76=4194\", \"381\u00d77=2667\"],\n  [\"128\u00d77=896\", \"182\u00d74=728\"],\n  [\"792\u00d74=3168\", \"382\u00d74=1528\"],\n  [\"878\u00d72=1756\", \"307\u00d77=2149\"],\n  [\"860\u00d77=6020\", \"859\u00d78=6872\"],\n  [\"161\u00d74=644\", \"715\u00d74=2860\"],\n  [\"880\u00d75=4400\", \"304\u00d79=2736\"],\n  [\"328\u00d76=1968\", \"379\u00d73=1137\"],\n  [\"289\u00d72=578\", \"507\u00d75=2535\"],\n  [\"520\u00d72=1040\", \"660\u00d73=1980\"],\n  [\"484\u00d79=4356\", \"288\u00d72=576\"],\n  [\"503\u00d75=2515\", \"299\u00d75=1495\"],\n  [\"983\u00d77=6881\", \"740\u00d72=1480\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication equations in the\n# table with their updated values, preserving all run/paragraph formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('364\u00d77=2548', '334\u00d76=2004'),\n    @('171\u00d73=513', '291\u00d76=1746'),\n    @('748\u00d79=6732', '239\u00d73=717'),\n    @('716\u00d77=5012', '945\u00d77=6615'),\n    @('318\u00d75=1590', '376\u00d77=2632'),\n    @('169\u00d72=338', '831\u00d72=1662'),\n    @('675\u00d75=3375', '713\u00d76=4278'),\n    @('201\u00d77=1407', '239\u00d73=717'),\n    @('568\u00d75=2840', '536\u00d76=3216'),\n    @('746\u00d72=1492', '464\u00d77=3248'),\n    @('520\u00d78=4160', '502\u00d77=3514'),\n    @('588\u00d78=4704', '401\u00d78=3208'),\n    @('699\u00d76=4194', '381\u00d77=2667'),\n    @('128\u00d77=896', '182\u00d74=728'),\n    @('792\u00d74=3168', '382\u00d74=1528'),\n    @('878\u00d72=1756', '307\u00d77=2149'),\n    @('860\u00d77=6020', '859\u00d78=6872'),\n    @('161\u00d74=644', '715\u00d74=2860'),\n    @('880\u00d75=4400', '304\u00d79=2736'),\n    @('328\u00d76=1968', '379\u00d73=1137'),\n    @('289\u00d72=578', '507\u00d75=2535'),\n    @('520\u00d72=1040', '660\u00d73=1980'),\n    @('484\u00d79=4356', '288\u00d72=576'),\n    @('503\u00d75=2515', '299\u00d75=1495'),\n    @('983\u00d77=6881', '740\u00d72=1480'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
